{"js": "// Helper: wrap a <w:p>...</w:p> (or other body-level) fragment into the\n// \"flat OPC\" package format Word.Range/Paragraph.insertOoxml expects.\nfunction wrapBodyOoxml(bodyFragmentXml) {\n  return '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n    'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' + bodyFragmentXml + '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>';\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the two paragraphs we need to touch by their (stable) text content\n// rather than a hard-coded index, so the script is resilient to minor\n// structural differences.\nlet fourthIndex = -1;\nlet mergeFieldIndex = -1;\nlet twentiethIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t === \"\u0427\u0435\u0442\u0432\u0435\u0440\u0442\u044b\u0439 \u0430\u0431\u0437\u0430\u0446\") {\n    fourthIndex = i;\n  } else if (t === \"\u00ab//@block1789\u00bb\") {\n    mergeFieldIndex = i;\n  } else if (t === \"\u0414\u0432\u0430\u0434\u0446\u0430\u0442\u044b\u0439 \u0430\u0431\u0437\u0430\u0446\") {\n    twentiethIndex = i;\n  }\n}\n\n// 1) The \"\u0427\u0435\u0442\u0432\u0435\u0440\u0442\u044b\u0439 \u0430\u0431\u0437\u0430\u0446\" paragraph currently carries the _GoBack\n//    bookmark at its end; rewrite it without the bookmark (the bookmark\n//    is relocated into the \"\u0414\u0432\u0430\u0434\u0446\u0430\u0442\u044b\u0439 \u0430\u0431\u0437\u0430\u0446\" paragraph below).\nif (fourthIndex !== -1) {\n  const p = paragraphs.items[fourthIndex];\n  const xml = \"<w:p><w:r><w:t>\u0427\u0435\u0442\u0432\u0435\u0440\u0442\u044b\u0439 \u0430\u0431\u0437\u0430\u0446</w:t></w:r></w:p>\";\n  p.insertOoxml(wrapBodyOoxml(xml), Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Replace the MERGEFIELD (\u00ab//@block1789\u00bb) paragraph with a plain\n//    paragraph containing the literal text \"//@block1789\" tagged en-US.\nif (mergeFieldIndex !== -1) {\n  const p = paragraphs.items[mergeFieldIndex];\n  const xml =\n    \"<w:p><w:pPr><w:rPr><w:lang w:val=\\\"en-US\\\"/></w:rPr></w:pPr>\" +\n    \"<w:r><w:rPr><w:lang w:val=\\\"en-US\\\"/></w:rPr><w:t>//@block1789</w:t></w:r></w:p>\";\n  p.insertOoxml(wrapBodyOoxml(xml), Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 3) Re-query paragraphs (content above may have shifted the collection)\n//    and split \"\u0414\u0432\u0430\u0434\u0446\u0430\u0442\u044b\u0439 \u0430\u0431\u0437\u0430\u0446\" into \"\u0414\u0432\u0430\u0434\u0446\u0430\u0442\u044b\u0439 \u0430\u0431\" + _GoBack + \"\u0437\u0430\u0446\".\nconst paragraphs2 = body.paragraphs;\nparagraphs2.load(\"items/text\");\nawait context.sync();\n\nlet targetIndex = -1;\nfor (let i = 0; i < paragraphs2.items.length; i++) {\n  if (paragraphs2.items[i].text === \"\u0414\u0432\u0430\u0434\u0446\u0430\u0442\u044b\u0439 \u0430\u0431\u0437\u0430\u0446\") {\n    targetIndex = i;\n    break;\n  }\n}\n\nif (targetIndex !== -1) {\n  const p = paragraphs2.items[targetIndex];\n  const xml =\n    \"<w:p><w:r><w:t>\u0414\u0432\u0430\u0434\u0446\u0430\u0442\u044b\u0439 \u0430\u0431</w:t></w:r>\" +\n    \"<w:bookmarkStart w:id=\\\"0\\\" w:name=\\\"_GoBack\\\"/><w:bookmarkEnd w:id=\\\"0\\\"/>\" +\n    \"<w:r><w:t>\u0437\u0430\u0446</w:t></w:r></w:p>\";\n  p.insertOoxml(wrapBodyOoxml(xml), Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Wrap-BodyOoxml($bodyFragmentXml) {\n@\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>$bodyFragmentXml</w:body></w:document>\n</pkg:xmlData></pkg:part></pkg:package>\n\"@\n}\n\n$fieldParaText = [char]0x00AB + \"//@block1789\" + [char]0x00BB + \"`r\"\n\n# 1) Remove the _GoBack bookmark from its current location (end of the\n#    \"\u0427\u0435\u0442\u0432\u0435\u0440\u0442\u044b\u0439 \u0430\u0431\u0437\u0430\u0446\" paragraph) - it is relocated below, into the\n#    \"\u0414\u0432\u0430\u0434\u0446\u0430\u0442\u044b\u0439 \u0430\u0431\u0437\u0430\u0446\" paragraph.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 2) Replace the MERGEFIELD paragraph (displayed as \u00ab//@block1789\u00bb) with a\n#    plain paragraph containing the literal text //@block1789, tagged en-US.\n#    Delete the field first so the paragraph collapses down to just its\n#    paragraph mark, then InsertXML the whole (mark-inclusive) range so the\n#    fresh paragraph's own pPr/rPr (the en-US language tag) survive.\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -eq $fieldParaText) {\n        $pRange = $p.Range\n        for ($fi = 1; $fi -le $d.Fields.Count; $fi++) {\n            $fld = $d.Fields.Item($fi)\n            if ($fld.Result.Start -ge $pRange.Start -and $fld.Result.End -le $pRange.End) {\n                $fld.Delete()\n                break\n            }\n        }\n        $p2 = $d.Paragraphs.Item($i)\n        $xml = \"<w:p><w:pPr><w:rPr><w:lang w:val=`\"en-US`\"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=`\"en-US`\"/></w:rPr><w:t>//@block1789</w:t></w:r></w:p>\"\n        $p2.Range.InsertXML((Wrap-BodyOoxml $xml))\n        break\n    }\n}\n\n# 3) Split \"\u0414\u0432\u0430\u0434\u0446\u0430\u0442\u044b\u0439 \u0430\u0431\u0437\u0430\u0446\" into \"\u0414\u0432\u0430\u0434\u0446\u0430\u0442\u044b\u0439 \u0430\u0431\" + _GoBack bookmark + \"\u0437\u0430\u0446\".\n#    This paragraph keeps real text, so InsertXML must target the range\n#    *excluding* the trailing paragraph mark, otherwise an extra empty\n#    paragraph gets minted.\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -eq \"\u0414\u0432\u0430\u0434\u0446\u0430\u0442\u044b\u0439 \u0430\u0431\u0437\u0430\u0446`r\") {\n        $r = $p.Range\n        $rNoMark = $d.Range($r.Start, $r.End - 1)\n        $xml = \"<w:p><w:r><w:t>\u0414\u0432\u0430\u0434\u0446\u0430\u0442\u044b\u0439 \u0430\u0431</w:t></w:r><w:bookmarkStart w:id=`\"0`\" w:name=`\"_GoBack`\"/><w:bookmarkEnd w:id=`\"0`\"/><w:r><w:t>\u0437\u0430\u0446</w:t></w:r></w:p>\"\n        $rNoMark.InsertXML((Wrap-BodyOoxml $xml))\n        break\n    }\n}\n"}
